$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, shifting existing rows 44-69 down to 45-70
$ws.Rows(44).Insert()

# Populate the newly inserted row 44 with this week's new data entry
$ws.Range("A44").Value = 9
$ws.Range("B44").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C44").Value = "Metropolitana"
$ws.Range("D44").Value = 44762
$ws.Range("E44").Value = 13
$ws.Range("F44").Value = 100112029
$ws.Range("G44").Value = "Orégano"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 16
$ws.Range("K44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("M44").Value = 20000
$ws.Range("N44").Value = "$/docena de atados"
$ws.Range("O44").Value = "Región Metropolitana"
$ws.Range("P44").Value = 6667
$ws.Range("Q44").Value = 3
$ws.Range("R44").Value = "Hortaliza"
